$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("A1:K9")

# Force text storage for the whole target range so numeric-looking
# strings (e.g. "12", "171.42") are kept as text, matching the source
# data (every cell in the sheet is authored as a string).
$target.NumberFormat = "@"

$ws.Cells.Item(1, 1).Value = "venue"
$ws.Cells.Item(1, 2).Value = "date"
$ws.Cells.Item(1, 3).Value = "result"
$ws.Cells.Item(1, 4).Value = "ownTeam"
$ws.Cells.Item(1, 5).Value = "oppTeam"
$ws.Cells.Item(1, 6).Value = "batsman"
$ws.Cells.Item(1, 7).Value = "totalRuns"
$ws.Cells.Item(1, 8).Value = "totalBalls"
$ws.Cells.Item(1, 9).Value = "total4s"
$ws.Cells.Item(1, 10).Value = "total6s"
$ws.Cells.Item(1, 11).Value = "sr"

$ws.Cells.Item(2, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(2, 2).Value = " October 27 2020"
$ws.Cells.Item(2, 3).Value = "Sunrisers won by 88 runs"
$ws.Cells.Item(2, 4).Value = "Delhi Capitals"
$ws.Cells.Item(2, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(2, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(2, 7).Value = "3"
$ws.Cells.Item(2, 8).Value = "7"
$ws.Cells.Item(2, 9).Value = "0"
$ws.Cells.Item(2, 10).Value = "0"
$ws.Cells.Item(2, 11).Value = "42.85"

$ws.Cells.Item(3, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(3, 2).Value = " November 05 2020"
$ws.Cells.Item(3, 3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(3, 4).Value = "Delhi Capitals"
$ws.Cells.Item(3, 5).Value = "Mumbai Indians"
$ws.Cells.Item(3, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(3, 7).Value = "15"
$ws.Cells.Item(3, 8).Value = "15"
$ws.Cells.Item(3, 9).Value = "2"
$ws.Cells.Item(3, 10).Value = "0"
$ws.Cells.Item(3, 11).Value = "100.00"

$ws.Cells.Item(4, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(4, 2).Value = " November 10 2020"
$ws.Cells.Item(4, 3).Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Cells.Item(4, 4).Value = "Delhi Capitals"
$ws.Cells.Item(4, 5).Value = "Mumbai Indians"
$ws.Cells.Item(4, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(4, 7).Value = "0"
$ws.Cells.Item(4, 8).Value = "0"
$ws.Cells.Item(4, 9).Value = "0"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).Value = "-"

$ws.Cells.Item(5, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(5, 2).Value = " October 31 2020"
$ws.Cells.Item(5, 3).Value = "Mumbai won by 9 wickets (with 34 balls remaining)"
$ws.Cells.Item(5, 4).Value = "Delhi Capitals"
$ws.Cells.Item(5, 5).Value = "Mumbai Indians"
$ws.Cells.Item(5, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(5, 7).Value = "12"
$ws.Cells.Item(5, 8).Value = "7"
$ws.Cells.Item(5, 9).Value = "0"
$ws.Cells.Item(5, 10).Value = "1"
$ws.Cells.Item(5, 11).Value = "171.42"

$ws.Cells.Item(6, 1).Value = " Sharjah"
$ws.Cells.Item(6, 2).Value = " October 09 2020"
$ws.Cells.Item(6, 3).Value = "Capitals won by 46 runs"
$ws.Cells.Item(6, 4).Value = "Delhi Capitals"
$ws.Cells.Item(6, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(6, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(6, 7).Value = "2"
$ws.Cells.Item(6, 8).Value = "3"
$ws.Cells.Item(6, 9).Value = "0"
$ws.Cells.Item(6, 10).Value = "0"
$ws.Cells.Item(6, 11).Value = "66.66"

$ws.Cells.Item(7, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(7, 2).Value = " September 20 2020"
$ws.Cells.Item(7, 3).Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Cells.Item(7, 4).Value = "Delhi Capitals"
$ws.Cells.Item(7, 5).Value = "Kings XI Punjab"
$ws.Cells.Item(7, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(7, 7).Value = "0"
$ws.Cells.Item(7, 8).Value = "0"
$ws.Cells.Item(7, 9).Value = "0"
$ws.Cells.Item(7, 10).Value = "0"
$ws.Cells.Item(7, 11).Value = "-"

$ws.Cells.Item(8, 1).Value = " Abu Dhabi"
$ws.Cells.Item(8, 2).Value = " October 24 2020"
$ws.Cells.Item(8, 3).Value = "KKR won by 59 runs"
$ws.Cells.Item(8, 4).Value = "Delhi Capitals"
$ws.Cells.Item(8, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(8, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(8, 7).Value = "9"
$ws.Cells.Item(8, 8).Value = "10"
$ws.Cells.Item(8, 9).Value = "1"
$ws.Cells.Item(8, 10).Value = "0"
$ws.Cells.Item(8, 11).Value = "90.00"

$ws.Cells.Item(9, 1).Value = " Abu Dhabi"
$ws.Cells.Item(9, 2).Value = " September 29 2020"
$ws.Cells.Item(9, 3).Value = "Sunrisers won by 15 runs"
$ws.Cells.Item(9, 4).Value = "Delhi Capitals"
$ws.Cells.Item(9, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(9, 6).Value = "Kagiso Rabada "
$ws.Cells.Item(9, 7).Value = "15"
$ws.Cells.Item(9, 8).Value = "7"
$ws.Cells.Item(9, 9).Value = "1"
$ws.Cells.Item(9, 10).Value = "1"
$ws.Cells.Item(9, 11).Value = "214.28"

# Reset formatting back to the workbook default style (no explicit
# NumberFormat override lingering on the cells).
$target.Style = "Normal"
